$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.643.87'
$ws.Cells.Item(2, 5).Value = '  +1.03%  '

$ws.Cells.Item(3, 4).Value = '2.437.48'
$ws.Cells.Item(3, 5).Value = '  +1.22%  '

$ws.Cells.Item(5, 4).Value = "'567.12"
$ws.Cells.Item(5, 5).Value = '  +0.78%  '

$ws.Cells.Item(6, 4).Value = "'145.24"
$ws.Cells.Item(6, 5).Value = '  +2.04%  '

$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 4).Value = "'0.534"
$ws.Cells.Item(8, 5).Value = '  +1.29%  '

$ws.Cells.Item(9, 5).Value = '  +1.99%  '

$ws.Cells.Item(10, 5).Value = '  +0.52%  '

$ws.Cells.Item(11, 5).Value = '  -0.66%  '

$ws.Cells.Item(12, 5).Value = '  +1.35%  '

$ws.Cells.Item(13, 4).Value = "'26.87"
$ws.Cells.Item(13, 5).Value = '  +5.31%  '

$ws.Cells.Item(14, 5).Value = '  +6.56%  '

$ws.Cells.Item(15, 4).Value = '2.877.58'
$ws.Cells.Item(15, 5).Value = '  +1.22%  '

$ws.Cells.Item(16, 4).Value = '62.404.34'
$ws.Cells.Item(16, 5).Value = '  +0.49%  '

$ws.Cells.Item(17, 4).Value = '2.440.50'
$ws.Cells.Item(17, 5).Value = '  +0.72%  '

$ws.Cells.Item(18, 4).Value = "'11.23"
$ws.Cells.Item(18, 5).Value = '  -0.42%  '

$ws.Cells.Item(19, 4).Value = "'6.93"
$ws.Cells.Item(19, 5).Value = '  +1.29%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).Value = "'323.84"
$ws.Cells.Item(20, 5).Value = '  +0.89%  '

$ws.Cells.Item(21, 2).Value = 'Polkadot'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(21, 4).Value = "'4.18"
$ws.Cells.Item(21, 5).Value = '  +1.01%  '

$ws.Cells.Item(22, 5).Value = '  -0.03%  '

$ws.Cells.Item(23, 4).Value = "'67.24"
$ws.Cells.Item(23, 5).Value = '  +2.09%  '

$ws.Cells.Item(24, 5).Value = '  +2.63%  '

$ws.Cells.Item(25, 4).Value = "'8.74"
$ws.Cells.Item(25, 5).Value = '  -0.64%  '

$ws.Cells.Item(26, 5).Value = '  +8.52%  '

$ws.Cells.Item(27, 4).Value = "'564.35"
$ws.Cells.Item(27, 5).Value = '  -0.25%  '

$ws.Cells.Item(28, 4).Value = '2.553.63'
$ws.Cells.Item(28, 5).Value = '  +1.01%  '

$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 5).Value = '  -0.15%  '

$ws.Cells.Item(30, 4).Value = "'8.39"
$ws.Cells.Item(30, 5).Value = '  +2.70%  '

$ws.Cells.Item(31, 4).Value = "'1.46"
$ws.Cells.Item(31, 5).Value = '  +3.49%  '

$ws.Cells.Item(32, 4).Value = "'0.147"
$ws.Cells.Item(32, 5).Value = '  -0.14%  '

$ws.Cells.Item(33, 5).Value = '  +0.17%  '

$ws.Cells.Item(34, 4).Value = "'1.54"
$ws.Cells.Item(34, 5).Value = '  +1.34%  '

$ws.Cells.Item(35, 4).Value = "'4.87"
$ws.Cells.Item(35, 5).Value = '  +4.15%  '

$ws.Cells.Item(37, 5).Value = '  +1.26%  '

$ws.Cells.Item(38, 4).Value = "'5.43"
$ws.Cells.Item(38, 5).Value = '  -0.36%  '

$ws.Cells.Item(39, 5).Value = '  +0.96%  '

$ws.Cells.Item(40, 4).Value = "'148.45"
$ws.Cells.Item(40, 5).Value = '  -2.26%  '

$ws.Cells.Item(41, 4).Value = "'1.82"
$ws.Cells.Item(41, 5).Value = '  +2.20%  '

$ws.Cells.Item(42, 5).Value = '  +0.18%  '

$ws.Cells.Item(43, 5).Value = '  +6.59%  '

$ws.Cells.Item(44, 4).Value = "'148.64"
$ws.Cells.Item(44, 5).Value = '  +0.73%  '

$ws.Cells.Item(45, 4).Value = "'3.68"
$ws.Cells.Item(45, 5).Value = '  +1.63%  '

$ws.Cells.Item(46, 5).Value = '  +1.21%  '

$ws.Cells.Item(47, 4).Value = "'20.53"
$ws.Cells.Item(47, 5).Value = '  +3.60%  '

$ws.Cells.Item(48, 4).Value = "'0.600"
$ws.Cells.Item(48, 5).Value = '  +1.69%  '

$ws.Cells.Item(49, 4).Value = "'0.0231"
$ws.Cells.Item(49, 5).Value = '  +3.05%  '

$ws.Cells.Item(50, 4).Value = "'0.0927"
$ws.Cells.Item(50, 5).Value = '  +1.38%  '

$ws.Cells.Item(51, 5).Value = '  +0.59%  '
